$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 79 (shifts existing rows 79-115 down to 80-116)
$ws.Rows.Item(79).EntireRow.Insert()

# Populate the newly inserted row 79 with this week's price record
# (same market/product metadata as the surrounding rows, new date + volume/price figures)
$ws.Cells.Item(79, 1).Value = 10
$ws.Cells.Item(79, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(79, 3).Value = "La Araucanía"
$ws.Cells.Item(79, 4).Value = 45146
$ws.Cells.Item(79, 5).Value = 9
$ws.Cells.Item(79, 6).Value = "Fruta"
$ws.Cells.Item(79, 7).Value = 100108
$ws.Cells.Item(79, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(79, 9).Value = 100108007
$ws.Cells.Item(79, 10).Value = "Coco"
$ws.Cells.Item(79, 11).Value = "Sin especificar"
$ws.Cells.Item(79, 12).Value = "Primera"
$ws.Cells.Item(79, 13).Value = 20
$ws.Cells.Item(79, 14).Value = 32000
$ws.Cells.Item(79, 15).Value = 32000
$ws.Cells.Item(79, 16).Value = 32000
$ws.Cells.Item(79, 17).Value = "$/malla 20 unidades"
$ws.Cells.Item(79, 18).Value = "Perú"
$ws.Cells.Item(79, 19).Value = 1600
$ws.Cells.Item(79, 20).Value = 20
